$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.211.90"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").Value = "2.056.58"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.42"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.617"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.68"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +7.95%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0831"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.02%  "
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.80"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.22%  "
$ws.Range("D13").Value = "2.362.44"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.05"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.760"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.31"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.93%  "
$ws.Range("D17").Value = "2.056.20"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "38.125.92"
$ws.Range("E18").Value = "  +2.40%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.78"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.19%  "
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "225.10"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.13"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.25"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.58%  "
$ws.Range("E28").Value = "  -0.73%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.98"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.29"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.94%  "
$ws.Range("E31").Value = "  +2.78%  "
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.55"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.04"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0606"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.08%  "
$ws.Range("E36").Value = "  +11.30%  "
$ws.Range("E37").Value = "  -2.24%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").Value = "1.527.43"
$ws.Range("E40").Value = "  +3.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "98.11"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("E42").Value = "  +1.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.79"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0929"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +0.56%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.01"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -7.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.99"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.89%  "
$ws.Range("E50").Value = "  -0.80%  "
$ws.Range("D51").Value = "2.250.69"
$ws.Range("E51").Value = "  +1.45%  "
